# ---------------------------------------------------------------------------
# Updates the team-specific state transition probability matrix
# (Grand Canyon_B) on Sheet1 with freshly computed empirical frequencies.
# This refreshes probabilities for rows 2-4,6-13,15-19 (row indices correspond
# to the "Starting_State" column in column A); the underlying logic that
# consumes this matrix in the simulation has not been wired up yet.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1890243902439024
$ws.Range("C2").Value = 0.551829268292683
$ws.Range("J2").Value = 0.01829268292682927
$ws.Range("P2").Value = 0.1310975609756098
$ws.Range("S2").Value = 0.1097560975609756

# Row 3
$ws.Range("B3").Value = 0.0374331550802139
$ws.Range("C3").Value = 0.0213903743315508
$ws.Range("J3").Value = 0.03208556149732621
$ws.Range("P3").Value = 0.7005347593582888
$ws.Range("S3").Value = 0.2085561497326203

# Row 4
$ws.Range("J4").Value = 0.02325581395348837
$ws.Range("P4").Value = 0.6744186046511628
$ws.Range("S4").Value = 0.3023255813953488

# Row 6
$ws.Range("B6").Value = 0.05063291139240506
$ws.Range("D6").Value = 0.008438818565400843
$ws.Range("F6").Value = 0.07172995780590717
$ws.Range("J6").Value = 0.2489451476793249
$ws.Range("O6").Value = 0.02109704641350211
$ws.Range("Q6").Value = 0.1518987341772152
$ws.Range("R6").Value = 0.05485232067510549
$ws.Range("S6").Value = 0.3924050632911392

# Row 7
$ws.Range("B7").Value = 0.08365019011406843
$ws.Range("D7").Value = 0.007604562737642586
$ws.Range("F7").Value = 0.05703422053231939
$ws.Range("J7").Value = 0.1368821292775665
$ws.Range("O7").Value = 0.01520912547528517
$ws.Range("Q7").Value = 0.1749049429657795
$ws.Range("R7").Value = 0.07224334600760456
$ws.Range("S7").Value = 0.4524714828897338

# Row 8
$ws.Range("B8").Value = 0.1118012422360248
$ws.Range("D8").Value = 0.01863354037267081
$ws.Range("F8").Value = 0.08281573498964803
$ws.Range("J8").Value = 0.09730848861283643
$ws.Range("O8").Value = 0.010351966873706
$ws.Range("Q8").Value = 0.1801242236024845
$ws.Range("R8").Value = 0.07039337474120083
$ws.Range("S8").Value = 0.4285714285714285

# Row 9
$ws.Range("B9").Value = 0.1058201058201058
$ws.Range("D9").Value = 0.01587301587301587
$ws.Range("F9").Value = 0.07407407407407407
$ws.Range("J9").Value = 0.1111111111111111
$ws.Range("O9").Value = 0.005291005291005291
$ws.Range("Q9").Value = 0.1957671957671958
$ws.Range("R9").Value = 0.06349206349206349
$ws.Range("S9").Value = 0.4285714285714285

# Row 10
$ws.Range("B10").Value = 0.1178369652945924
$ws.Range("D10").Value = 0.02179176755447942
$ws.Range("E10").Value = 0.002421307506053269
$ws.Range("F10").Value = 0.06941081517352704
$ws.Range("J10").Value = 0.132364810330912
$ws.Range("O10").Value = 0.01856335754640839
$ws.Range("Q10").Value = 0.1767554479418886
$ws.Range("R10").Value = 0.07102502017756256
$ws.Range("S10").Value = 0.3898305084745763

# Row 11
$ws.Range("G11").Value = 0.1417721518987342
$ws.Range("J11").Value = 0.05569620253164557
$ws.Range("K11").Value = 0.1746835443037975
$ws.Range("L11").Value = 0.6025316455696202
$ws.Range("S11").Value = 0.02531645569620253

# Row 12
$ws.Range("G12").Value = 0.7625
$ws.Range("J12").Value = 0.1375
$ws.Range("K12").Value = 0.004166666666666667
$ws.Range("L12").Value = 0.025
$ws.Range("S12").Value = 0.07083333333333333

# Row 13
$ws.Range("G13").Value = 0.6875
$ws.Range("S13").Value = 0.0625

# Row 15
$ws.Range("F15").Value = 0.0091324200913242
$ws.Range("H15").Value = 0.1917808219178082
$ws.Range("I15").Value = 0.045662100456621
$ws.Range("J15").Value = 0.3561643835616438
$ws.Range("K15").Value = 0.0684931506849315
$ws.Range("M15").Value = 0.0136986301369863
$ws.Range("N15").Value = 0.0045662100456621
$ws.Range("O15").Value = 0.0365296803652968
$ws.Range("S15").Value = 0.273972602739726

# Row 16
$ws.Range("F16").Value = 0.02487562189054726
$ws.Range("H16").Value = 0.154228855721393
$ws.Range("I16").Value = 0.05970149253731343
$ws.Range("J16").Value = 0.3830845771144278
$ws.Range("K16").Value = 0.1293532338308458
$ws.Range("M16").Value = 0.02487562189054726
$ws.Range("N16").Value = 0.004975124378109453
$ws.Range("O16").Value = 0.07462686567164178
$ws.Range("S16").Value = 0.1442786069651741

# Row 17
$ws.Range("F17").Value = 0.01650943396226415
$ws.Range("H17").Value = 0.1509433962264151
$ws.Range("I17").Value = 0.08490566037735849
$ws.Range("J17").Value = 0.4150943396226415
$ws.Range("K17").Value = 0.1202830188679245
$ws.Range("M17").Value = 0.01650943396226415
$ws.Range("O17").Value = 0.05660377358490566
$ws.Range("S17").Value = 0.1391509433962264

# Row 18
$ws.Range("F18").Value = 0.01851851851851852
$ws.Range("H18").Value = 0.1666666666666667
$ws.Range("I18").Value = 0.1111111111111111
$ws.Range("J18").Value = 0.3580246913580247
$ws.Range("K18").Value = 0.1234567901234568
$ws.Range("M18").Value = 0.02469135802469136
$ws.Range("O18").Value = 0.05555555555555555
$ws.Range("S18").Value = 0.1419753086419753

# Row 19
$ws.Range("F19").Value = 0.01683029453015428
$ws.Range("H19").Value = 0.2279102384291725
$ws.Range("I19").Value = 0.08064516129032258
$ws.Range("J19").Value = 0.3253856942496494
$ws.Range("K19").Value = 0.1430575035063114
$ws.Range("M19").Value = 0.02103786816269285
$ws.Range("N19").Value = 0.002103786816269285
$ws.Range("O19").Value = 0.06661991584852735
$ws.Range("S19").Value = 0.1164095371669004

